# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, shifting the existing "Late" / "heading" / "Outstanding" columns
# one place to the right, then leave that sheet active/selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Capture column M's width (character units) so the freshly inserted
# column N can match it exactly, then insert the new blank column.
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab with K21 selected
# (this also clears tabSelected on whichever sheet was active before).
$ws.Activate()
$ws.Range("K21").Select() | Out-Null
